$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new rows (22-30) continuing the existing data pattern:
# regcntr_id 10002..10010, machine_id 10021..10029, lang_code "eng",
# is_active TRUE, cr_by "superadmin", cr_dtimes/eff_dtimes "now()"
for ($i = 0; $i -lt 9; $i++) {
    $row = 22 + $i
    $ws.Cells.Item($row, 1).Value = 10002 + $i
    $ws.Cells.Item($row, 2).Value = 10021 + $i
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Match the saved view state: scrolled down with B22:B30 selected
$ws.Range("A19").Select()
$ws.Range("B22:B30").Select()

# Page setup: portrait orientation (as added in pageSetup element)
$ws.PageSetup.Orientation = 1
